$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values stay text (avoid Excel auto-converting numeric-looking strings to numbers)
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '25.589.99'
$ws.Range('E2').Value = '  -1.63%  '
$ws.Range('D3').Value = '1.590.61'
$ws.Range('E3').Value = '  -2.89%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '208.91'
$ws.Range('E5').Value = '  -2.74%  '
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('E7').Value = '  -4.72%  '
$ws.Range('E8').Value = '  -2.17%  '
$ws.Range('E9').Value = '  -2.12%  '
$ws.Range('E10').Value = '  -2.93%  '
$ws.Range('D11').Value = '0.0785'
$ws.Range('E11').Value = '  -0.65%  '
$ws.Range('D12').Value = '1.811.17'
$ws.Range('E12').Value = '  -2.93%  '
$ws.Range('D13').Value = '1.591.17'
$ws.Range('E13').Value = '  -2.82%  '
$ws.Range('E14').Value = '  -3.93%  '
$ws.Range('E15').Value = '  -3.80%  '
$ws.Range('D16').Value = '25.582.24'
$ws.Range('E16').Value = '  -1.68%  '
$ws.Range('E17').Value = '  -2.31%  '
$ws.Range('D18').Value = '0.0₃0712'
$ws.Range('E18').Value = '  -4.30%  '
$ws.Range('E19').Value = '  +0.02%  '
$ws.Range('D20').Value = '188.49'
$ws.Range('E20').Value = '  -1.83%  '
$ws.Range('E21').Value = '  -1.89%  '
$ws.Range('E22').Value = '  -4.03%  '
$ws.Range('D23').Value = '5.94'
$ws.Range('E23').Value = '  -2.67%  '
$ws.Range('E24').Value = '  -0.02%  '
$ws.Range('E25').Value = '  -4.21%  '
$ws.Range('D26').Value = '140.69'
$ws.Range('E26').Value = '  -2.27%  '
$ws.Range('E27').Value = '  -5.03%  '
$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D28').Value = '6.50'
$ws.Range('E28').Value = '  -4.88%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Value = '14.97'
$ws.Range('E29').Value = '  -1.76%  '
$ws.Range('E30').Value = '  -4.40%  '
$ws.Range('D31').Value = '0.0466'
$ws.Range('E31').Value = '  -3.81%  '
$ws.Range('E32').Value = '  -2.57%  '
$ws.Range('D33').Value = '3.01'
$ws.Range('E33').Value = '  -4.29%  '
$ws.Range('D34').Value = '2.39'
$ws.Range('E34').Value = '  -0.82%  '
$ws.Range('E35').Value = '  -1.68%  '
$ws.Range('D36').Value = '1.090.95'
$ws.Range('E36').Value = '  -4.03%  '
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('E38').Value = '  -3.82%  '
$ws.Range('E39').Value = '  -2.68%  '
$ws.Range('D40').Value = '0.781'
$ws.Range('E40').Value = '  -10.06%  '
$ws.Range('E41').Value = '  -4.62%  '
$ws.Range('D42').Value = '95.26'
$ws.Range('E42').Value = '  -3.29%  '
$ws.Range('D43').Value = '1.724.78'
$ws.Range('E43').Value = '  -2.87%  '
$ws.Range('D44').Value = '5.07'
$ws.Range('E44').Value = '  -3.19%  '
$ws.Range('D45').Value = '0.735'
$ws.Range('E45').Value = '  -5.74%  '
$ws.Range('D46').Value = '0.0₆0107'
$ws.Range('E46').Value = '  -6.98%  '
$ws.Range('D47').Value = '52.99'
$ws.Range('E47').Value = '  -3.86%  '
$ws.Range('E48').Value = '  -3.77%  '
$ws.Range('D49').Value = '1.43'
$ws.Range('E49').Value = '  -4.21%  '
$ws.Range('E50').Value = '  -1.31%  '
$ws.Range('E51').Value = '  -0.07%  '

# Restore default (styleless) formatting for column D so XML matches original (no explicit style index)
$ws.Range('D2:D51').Style = 'Normal'
